$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gameflow")

# Update the dialogue line for the Tanuki's second line to include the
# screen-shake text event tag.
$ws.Range("C4").Value = "Let's [screen-shake=0.2,0.3]get em!"

# Update the active selection to match the authored change.
$ws.Range("C4").Select()
